# Update the "Price" (D) and "Volume(1h)" (E) columns of the cryptos list
# with refreshed values from the scraper. Numeric-looking price strings are
# entered with a leading apostrophe so Excel keeps them as text (preserving
# trailing zeros / exact formatting) instead of converting them to numbers.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "24.613.89"
$ws.Range("E2").Value = "  +1.53%  "
$ws.Range("D3").Value = "1.703.92"
$ws.Range("E3").Value = "  +1.34%  "
$ws.Range("D4").Value = "'1.003"
$ws.Range("E4").Value = "  +0.06%  "
$ws.Range("D5").Value = "'308.08"
$ws.Range("E5").Value = "  -0.57%  "
$ws.Range("D6").Value = "'0.9971"
$ws.Range("E6").Value = "  -0.05%  "
$ws.Range("D7").Value = "'0.3721"
$ws.Range("E7").Value = "  -0.68%  "
$ws.Range("D8").Value = "'48.95"
$ws.Range("E8").Value = "  +2.35%  "
$ws.Range("D9").Value = "'0.3433"
$ws.Range("E9").Value = "  -0.75%  "
$ws.Range("D10").Value = "'1.179"
$ws.Range("E10").Value = "  -1.69%  "
$ws.Range("D11").Value = "'0.07425"
$ws.Range("E11").Value = "  +1.13%  "
$ws.Range("D12").Value = "'0.9998"
$ws.Range("E12").Value = "  +0.02%  "
$ws.Range("D13").Value = "'20.80"
$ws.Range("E13").Value = "  +1.08%  "
$ws.Range("D14").Value = "'6.207"
$ws.Range("E14").Value = "  +1.22%  "
$ws.Range("D15").Value = "'6.914"
$ws.Range("E15").Value = "  +1.63%  "
$ws.Range("D16").Value = "1.706.08"
$ws.Range("E17").Value = "  +0.06%  "
$ws.Range("D18").Value = "'0.9972"
$ws.Range("E18").Value = "  -0.04%  "
$ws.Range("D19").Value = "'0.06688"
$ws.Range("E19").Value = "  -0.63%  "
$ws.Range("D20").Value = "'83.09"
$ws.Range("E20").Value = "  +1.08%  "
$ws.Range("D21").Value = "'17.01"
$ws.Range("E21").Value = "  +2.53%  "
$ws.Range("D22").Value = "'6.326"
$ws.Range("E22").Value = "  +2.78%  "
$ws.Range("D23").Value = "'13.09"
$ws.Range("E23").Value = "  +8.48%  "
$ws.Range("D24").Value = "24.656.52"
$ws.Range("E24").Value = "  +1.77%  "
$ws.Range("D25").Value = "'2.414"
$ws.Range("E25").Value = "  +0.15%  "
$ws.Range("D26").Value = "'2.758"
$ws.Range("E26").Value = "  +2.45%  "
$ws.Range("D27").Value = "'20.08"
$ws.Range("E27").Value = "  +2.09%  "
$ws.Range("D28").Value = "'149.25"
$ws.Range("E28").Value = "  -1.82%  "
$ws.Range("D29").Value = "'130.87"
$ws.Range("E29").Value = "  +2.70%  "
$ws.Range("D30").Value = "1.894.98"
$ws.Range("E30").Value = "  +1.71%  "
$ws.Range("D31").Value = "'1.169"
$ws.Range("E31").Value = "  +16.81%  "
$ws.Range("D32").Value = "'6.706"
$ws.Range("E32").Value = "  +2.47%  "
$ws.Range("D33").Value = "'4.196"
$ws.Range("E33").Value = "  +3.30%  "
$ws.Range("D34").Value = "'0.08770"
$ws.Range("E34").Value = "  +3.02%  "
$ws.Range("D35").Value = "'1.768"
$ws.Range("E35").Value = "  -0.93%  "
$ws.Range("D36").Value = "'13.57"
$ws.Range("E36").Value = "  +7.67%  "
$ws.Range("D37").Value = "'5.492"
$ws.Range("E37").Value = "  +1.38%  "
$ws.Range("D38").Value = "'0.06504"
$ws.Range("E38").Value = "  +0.23%  "
$ws.Range("D39").Value = "'8.882"
$ws.Range("E39").Value = "  -0.74%  "
$ws.Range("D42").Value = "'1.272"
$ws.Range("E42").Value = "  -0.98%  "
$ws.Range("D43").Value = "'0.6372"
$ws.Range("E43").Value = "  +2.28%  "
$ws.Range("D44").Value = "'0.9972"
$ws.Range("E44").Value = "  +0.00%  "
$ws.Range("D45").Value = "'13.82"
$ws.Range("E45").Value = "  +3.29%  "
$ws.Range("D46").Value = "'0.6053"
$ws.Range("E46").Value = "  +1.34%  "
$ws.Range("D47").Value = "'3.791"
$ws.Range("E47").Value = "  -0.69%  "
$ws.Range("D48").Value = "'2.103"
$ws.Range("E48").Value = "  +3.22%  "
$ws.Range("D49").Value = "'128.56"
$ws.Range("E49").Value = "  +0.36%  "
$ws.Range("D50").Value = "'0.07252"
$ws.Range("E50").Value = "  +1.05%  "
$ws.Range("D51").Value = "'78.75"

# Row 40/41 swap (Algorand/VeChain ranking order changed)
$ws.Range("B40").Value = "Algorand"
$ws.Range("C40").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D40").Value = "'0.2210"
$ws.Range("E40").Value = "  +3.39%  "
$ws.Range("B41").Value = "VeChain"
$ws.Range("C41").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D41").Value = "'0.02354"
$ws.Range("E41").Value = "  -0.20%  "
